$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text-like values (non-numeric strings) can be assigned directly.
$ws.Range("D2").Value = "41.631.17"
$ws.Range("E2").Value = "  +4.24%  "
$ws.Range("D3").Value = "2.255.65"
$ws.Range("E3").Value = "  +2.32%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("E5").Value = "  +3.01%  "
$ws.Range("E6").Value = "  +4.60%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +2.13%  "
$ws.Range("E10").Value = "  +4.87%  "
$ws.Range("E11").Value = "  +3.71%  "
$ws.Range("E12").Value = "  +2.06%  "
$ws.Range("E13").Value = "  +1.45%  "
$ws.Range("E14").Value = "  +2.82%  "
$ws.Range("D15").Value = "2.610.73"
$ws.Range("E15").Value = "  +2.63%  "
$ws.Range("E16").Value = "  +2.53%  "
$ws.Range("D17").Value = "2.259.69"
$ws.Range("E17").Value = "  +2.62%  "
$ws.Range("E18").Value = "  +3.25%  "
$ws.Range("D19").Value = "41.557.22"
$ws.Range("E19").Value = "  +4.26%  "
$ws.Range("E20").Value = "  +10.53%  "
$ws.Range("D21").Value = "0.0₃0900"
$ws.Range("E21").Value = "  +1.73%  "
$ws.Range("E22").Value = "  +2.43%  "
$ws.Range("E23").Value = "  +1.56%  "
$ws.Range("E24").Value = "  +2.08%  "
$ws.Range("E25").Value = "  +4.56%  "
$ws.Range("E27").Value = "  +5.24%  "
$ws.Range("E28").Value = "  +4.25%  "
$ws.Range("E29").Value = "  +2.06%  "
$ws.Range("E30").Value = "  -0.79%  "
$ws.Range("E31").Value = "  +0.53%  "
$ws.Range("E32").Value = "  +7.91%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("E34").Value = "  +3.74%  "
$ws.Range("E35").Value = "  +4.18%  "
$ws.Range("E36").Value = "  -0.71%  "
$ws.Range("E38").Value = "  +2.71%  "
$ws.Range("E39").Value = "  +6.76%  "
$ws.Range("E40").Value = "  +3.75%  "
$ws.Range("E41").Value = "  +2.62%  "
$ws.Range("E42").Value = "  +4.54%  "
$ws.Range("D43").Value = "2.059.41"
$ws.Range("E43").Value = "  -0.59%  "
$ws.Range("E44").Value = "  +1.10%  "
$ws.Range("E45").Value = "  +2.43%  "
$ws.Range("E46").Value = "  +2.51%  "
$ws.Range("E47").Value = "  +5.28%  "
$ws.Range("E48").Value = "  +3.83%  "
$ws.Range("E49").Value = "  +7.35%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("E50").Value = "  +3.83%  "
$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("E51").Value = "  +2.83%  "

# Numeric-looking strings must be forced to Text format to preserve
# the original inline-string representation (matching source data),
# then the number format is reset to Normal style to avoid leaving
# a visible text-format style applied to the cell.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "303.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "90.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.05"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.31"
$ws.Range("D11").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.12"
$ws.Range("D16").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.43"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "239.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.58"
$ws.Range("D25").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.46"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "160.10"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "34.16"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.50"
$ws.Range("D39").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.89"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.33"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.10"
$ws.Range("D46").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "72.22"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.50"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.14"
$ws.Range("D51").Style = "Normal"
